# refactor data gathering steps
# Adds "scraping_start_date" (column I) and "scraping_end_date" (column J)
# to the timeline sheet, mirroring the existing formation_*/testing_*
# columns (E/F/G/H) but holding their own scraping window values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

function Set-TimelineCell($row, $col, $value, $isHeader) {
    $cell = $ws.Cells.Item($row, $col)
    if (-not $isHeader) {
        $cell.NumberFormat = "@"
    }
    $cell.HorizontalAlignment = $xlCenter
    $cell.Value = $value
}

# J1 header ("scraping_end_date") is entered first.
Set-TimelineCell 1 10 "scraping_end_date" $true

# scraping_start_date values for the first block of rows (I2:I5).
Set-TimelineCell 2 9 "2019-07-25" $false
Set-TimelineCell 3 9 "2019-08-15" $false
Set-TimelineCell 4 9 "2019-10-14" $false
Set-TimelineCell 5 9 "2019-11-12" $false

# I1 header ("scraping_start_date").
Set-TimelineCell 1 9 "scraping_start_date" $true

# scraping_start_date values for the remaining rows (I6:I9).
Set-TimelineCell 6 9 "2020-12-29" $false
Set-TimelineCell 7 9 "2022-11-04" $false
Set-TimelineCell 8 9 "2022-12-06" $false
Set-TimelineCell 9 9 "2023-06-12" $false

# scraping_end_date values mirror the existing testing_end_date values,
# so they simply reuse the strings already present in column H.
$data = @(
    @{ Row = 2; End = "2019-09-07" },
    @{ Row = 3; End = "2019-09-28" },
    @{ Row = 4; End = "2019-11-27" },
    @{ Row = 5; End = "2019-12-28" },
    @{ Row = 6; End = "2021-02-16" },
    @{ Row = 7; End = "2022-12-20" },
    @{ Row = 8; End = "2023-01-19" },
    @{ Row = 9; End = "2023-08-01" }
)

foreach ($item in $data) {
    Set-TimelineCell $item.Row 10 $item.End $false
}

# Match column widths from the diff (col I: 16.7109375, col J: 18)
$ws.Columns.Item(9).ColumnWidth = 15.85
$ws.Columns.Item(10).ColumnWidth = 17.1

# Update the selected cell as in the diff (activeCell I9)
$ws.Range("I9").Select()
